$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 5987
$ws.Range("F8").Value = 10210
$ws.Range("F9").Value = 4037
$ws.Range("F14").Value = 716
$ws.Range("F15").Value = 4027
$ws.Range("F16").Value = 29
$ws.Range("F19").Value = 5694
$ws.Range("F21").Value = 2236
$ws.Range("F22").Value = 153
$ws.Range("F23").Value = 412
$ws.Range("F24").Value = 8543
$ws.Range("F26").Value = 95
$ws.Range("F28").Value = 2253
$ws.Range("F29").Value = 2301
$ws.Range("F30").Value = 1356
$ws.Range("F31").Value = 203
$ws.Range("F32").Value = 1881
$ws.Range("F34").Value = 307
$ws.Range("F38").Value = 33
$ws.Range("F44").Value = 1418
$ws.Range("F45").Value = 2292
$ws.Range("F46").Value = 171
$ws.Range("F47").Value = 260
$ws.Range("F49").Value = 15
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 158
$ws.Range("F12").Value = 138
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 643
$ws.Range("F4").Value = 83
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 643
$ws.Range("F6").Value = 5987
$ws.Range("F7").Value = 10210
$ws.Range("F8").Value = 4037
$ws.Range("F11").Value = 158
$ws.Range("F14").Value = 716
$ws.Range("F15").Value = 4027
$ws.Range("F16").Value = 29
$ws.Range("F19").Value = 5694
$ws.Range("F21").Value = 153
$ws.Range("F22").Value = 412
$ws.Range("F23").Value = 8543
$ws.Range("F24").Value = 138
$ws.Range("F26").Value = 95
$ws.Range("F28").Value = 2253
$ws.Range("F29").Value = 1356
$ws.Range("F30").Value = 203
$ws.Range("F31").Value = 1881
$ws.Range("F33").Value = 307
$ws.Range("F36").Value = 33
$ws.Range("F41").Value = 1418
$ws.Range("F43").Value = 2292
$ws.Range("F44").Value = 171
$ws.Range("F46").Value = 260
